$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1519
$ws.Range("F5").Value = 7802
$ws.Range("F6").Value = 4880
$ws.Range("F7").Value = 7168
$ws.Range("F8").Value = 294
$ws.Range("F9").Value = 1526
$ws.Range("F14").Value = 191
$ws.Range("F15").Value = 560
$ws.Range("F16").Value = 22
$ws.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202405/Q6qMlV7u1715918264358.png"
$ws.Range("F17").Value = 245
$ws.Range("F19").Value = 36
$ws.Range("F20").Value = 1230
$ws.Range("F24").Value = 1258
$ws.Range("F28").Value = 15
$ws.Range("F30").Value = 219
$ws.Range("F31").Value = 23
$ws.Range("F33").Value = 13
$ws.Range("F34").Value = 150
$ws.Range("F37").Value = 555
$ws.Range("F38").Value = 569
$ws.Range("F39").Value = 427
$ws.Range("F40").Value = 89
$ws.Range("F42").Value = 103
$ws.Range("F43").Value = 424

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 139
$ws.Range("F14").Value = 1736
$ws.Range("F15").Value = 565
$ws.Range("F17").Value = 20
$ws.Range("F29").Value = 27
$ws.Range("F32").Value = 877
$ws.Range("F34").Value = 999
$ws.Range("F41").Value = 108

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 738
$ws.Range("F6").Value = 699
$ws.Range("F7").Value = 174
$ws.Range("F8").Value = 104
$ws.Range("F9").Value = 1732
$ws.Range("F10").Value = 2643

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 738
$ws.Range("F3").Value = 1519
$ws.Range("F6").Value = 699
$ws.Range("F7").Value = 7802
$ws.Range("F8").Value = 174
$ws.Range("F9").Value = 4880
$ws.Range("F10").Value = 7168
$ws.Range("F11").Value = 294
$ws.Range("F12").Value = 1526
$ws.Range("F14").Value = 104
$ws.Range("F16").Value = 1732
$ws.Range("F17").Value = 2643
$ws.Range("F20").Value = 191
$ws.Range("F21").Value = 22
$ws.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202405/Q6qMlV7u1715918264358.png"
$ws.Range("F22").Value = 245
$ws.Range("F23").Value = 1230
$ws.Range("F26").Value = 1258
$ws.Range("F28").Value = 15
$ws.Range("F29").Value = 219
$ws.Range("F32").Value = 27
$ws.Range("F33").Value = 877
$ws.Range("F34").Value = 150
$ws.Range("F37").Value = 999
$ws.Range("F38").Value = 569
$ws.Range("F40").Value = 89
$ws.Range("F41").Value = 103
$ws.Range("F43").Value = 424
$ws.Range("F45").Value = 108
